$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("G2").Value = 33.64214033333334
$ws.Range("H2").Value = 100.926421
$ws.Range("I2").Value = 0.106995191696894
$ws.Range("J2").Value = 0.106995191696894
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.011299
$ws.Range("N2").Value = 0.033897
$ws.Range("O2").Value = 0.3524439315012944
$ws.Range("P2").Value = 0.3524439315012944
$ws.Range("Q2").Value = 0.3801225436263333
$ws.Range("R2").Value = 3.421102892637
$ws.Range("S2").Value = 0.03770980601338797
$ws.Range("T2").Value = 0.03770980601338797

$ws.Range("G3").Value = 33.64214033333334
$ws.Range("H3").Value = 100.926421
$ws.Range("I3").Value = 0.106995191696894
$ws.Range("J3").Value = 0.106995191696894
$ws.Range("O3").Value = 0.01059504871227008
$ws.Range("P3").Value = 0.01059504871227008
$ws.Range("Q3").Value = 0.01142711366655556
$ws.Range("R3").Value = 0.102844022999
$ws.Range("S3").Value = 0.001133619268007268
$ws.Range("T3").Value = 0.001133619268007267

$ws.Range("G4").Value = 33.64214033333334
$ws.Range("H4").Value = 100.926421
$ws.Range("I4").Value = 0.106995191696894
$ws.Range("J4").Value = 0.106995191696894
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02042033333333334
$ws.Range("N4").Value = 0.061261
$ws.Range("O4").Value = 0.6369610197864354
$ws.Range("P4").Value = 0.6369610197864355
$ws.Range("Q4").Value = 0.6869837196534446
$ws.Range("R4").Value = 6.182853476881
$ws.Range("S4").Value = 0.06815176641549875
$ws.Range("T4").Value = 0.06815176641549875

$ws.Range("H5").Value = 727.0751789999999
$ws.Range("I5").Value = 0.7707946777896593
$ws.Range("J5").Value = 0.7707946777896592
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.011299
$ws.Range("N5").Value = 0.033897
$ws.Range("O5").Value = 0.3524439315012944
$ws.Range("P5").Value = 0.3524439315012944
$ws.Range("Q5").Value = 2.738407482506999
$ws.Range("R5").Value = 24.64566734256299
$ws.Range("S5").Value = 0.271661906620461
$ws.Range("T5").Value = 0.271661906620461

$ws.Range("H6").Value = 727.0751789999999
$ws.Range("I6").Value = 0.7707946777896593
$ws.Range("J6").Value = 0.7707946777896592
$ws.Range("O6").Value = 0.01059504871227008
$ws.Range("P6").Value = 0.01059504871227008
$ws.Range("Q6").Value = 0.08232106748899999
$ws.Range("R6").Value = 0.7408896074009999
$ws.Range("S6").Value = 0.008166607158339964
$ws.Range("T6").Value = 0.008166607158339962

$ws.Range("H7").Value = 727.0751789999999
$ws.Range("I7").Value = 0.7707946777896593
$ws.Range("J7").Value = 0.7707946777896592
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02042033333333334
$ws.Range("N7").Value = 0.061261
$ws.Range("O7").Value = 0.6369610197864354
$ws.Range("P7").Value = 0.6369610197864355
$ws.Range("Q7").Value = 4.949039171191
$ws.Range("R7").Value = 44.541352540719
$ws.Range("S7").Value = 0.4909661640108582
$ws.Range("T7").Value = 0.4909661640108582

$ws.Range("G8").Value = 9.788214000000002
$ws.Range("H8").Value = 29.364642
$ws.Range("I8").Value = 0.0311303568359039
$ws.Range("J8").Value = 0.03113035683590389
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.011299
$ws.Range("N8").Value = 0.033897
$ws.Range("O8").Value = 0.3524439315012944
$ws.Range("P8").Value = 0.3524439315012944
$ws.Range("Q8").Value = 0.110597029986
$ws.Range("R8").Value = 0.995373269874
$ws.Range("S8").Value = 0.01097170535228417
$ws.Range("T8").Value = 0.01097170535228417

$ws.Range("G9").Value = 9.788214000000002
$ws.Range("H9").Value = 29.364642
$ws.Range("I9").Value = 0.0311303568359039
$ws.Range("J9").Value = 0.03113035683590389
$ws.Range("O9").Value = 0.01059504871227008
$ws.Range("P9").Value = 0.01059504871227008
$ws.Range("Q9").Value = 0.003324730022
$ws.Range("R9").Value = 0.029922570198
$ws.Range("S9").Value = 0.0003298276471067518
$ws.Range("T9").Value = 0.0003298276471067517

$ws.Range("G10").Value = 9.788214000000002
$ws.Range("H10").Value = 29.364642
$ws.Range("I10").Value = 0.0311303568359039
$ws.Range("J10").Value = 0.03113035683590389
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02042033333333334
$ws.Range("N10").Value = 0.061261
$ws.Range("O10").Value = 0.6369610197864354
$ws.Range("P10").Value = 0.6369610197864355
$ws.Range("Q10").Value = 0.1998785926180001
$ws.Range("R10").Value = 1.798907333562
$ws.Range("S10").Value = 0.01982882383651298
$ws.Range("T10").Value = 0.01982882383651298

$ws.Range("G11").Value = 12.12016933333333
$ws.Range("H11").Value = 36.360508
$ws.Range("I11").Value = 0.03854688876420623
$ws.Range("J11").Value = 0.03854688876420622
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.011299
$ws.Range("N11").Value = 0.033897
$ws.Range("O11").Value = 0.3524439315012944
$ws.Range("P11").Value = 0.3524439315012944
$ws.Range("Q11").Value = 0.1369457932973333
$ws.Range("R11").Value = 1.232512139676
$ws.Range("S11").Value = 0.01358561702319992
$ws.Range("T11").Value = 0.01358561702319991

$ws.Range("G12").Value = 12.12016933333333
$ws.Range("H12").Value = 36.360508
$ws.Range("I12").Value = 0.03854688876420623
$ws.Range("J12").Value = 0.03854688876420622
$ws.Range("O12").Value = 0.01059504871227008
$ws.Range("P12").Value = 0.01059504871227008
$ws.Range("Q12").Value = 0.004116817516888889
$ws.Range("R12").Value = 0.037051357652
$ws.Range("S12").Value = 0.0004084061641632214
$ws.Range("T12").Value = 0.0004084061641632213

$ws.Range("G13").Value = 12.12016933333333
$ws.Range("H13").Value = 36.360508
$ws.Range("I13").Value = 0.03854688876420623
$ws.Range("J13").Value = 0.03854688876420622
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.02042033333333334
$ws.Range("N13").Value = 0.061261
$ws.Range("O13").Value = 0.6369610197864354
$ws.Range("P13").Value = 0.6369610197864355
$ws.Range("Q13").Value = 0.2474978978431112
$ws.Range("R13").Value = 2.227481080588
$ws.Range("S13").Value = 0.02455286557684309
$ws.Range("T13").Value = 0.02455286557684309

$ws.Range("G14").Value = 16.51773933333333
$ws.Range("H14").Value = 49.553218
$ws.Range("I14").Value = 0.0525328849133368
$ws.Range("J14").Value = 0.05253288491333678
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.011299
$ws.Range("N14").Value = 0.033897
$ws.Range("O14").Value = 0.3524439315012944
$ws.Range("P14").Value = 0.3524439315012944
$ws.Range("Q14").Value = 0.1866339367273333
$ws.Range("R14").Value = 1.679705430546
$ws.Range("S14").Value = 0.01851489649196145
$ws.Range("T14").Value = 0.01851489649196145

$ws.Range("G15").Value = 16.51773933333333
$ws.Range("H15").Value = 49.553218
$ws.Range("I15").Value = 0.0525328849133368
$ws.Range("J15").Value = 0.05253288491333678
$ws.Range("O15").Value = 0.01059504871227008
$ws.Range("P15").Value = 0.01059504871227008
$ws.Range("Q15").Value = 0.005610525460222223
$ws.Range("R15").Value = 0.050494729142
$ws.Range("S15").Value = 0.0005565884746528816
$ws.Range("T15").Value = 0.0005565884746528815

$ws.Range("G16").Value = 16.51773933333333
$ws.Range("H16").Value = 49.553218
$ws.Range("I16").Value = 0.0525328849133368
$ws.Range("J16").Value = 0.05253288491333678
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02042033333333334
$ws.Range("N16").Value = 0.061261
$ws.Range("O16").Value = 0.6369610197864354
$ws.Range("P16").Value = 0.6369610197864355
$ws.Range("Q16").Value = 0.3372977430997778
$ws.Range("R16").Value = 3.035679687898
$ws.Range("S16").Value = 0.03346139994672245
$ws.Range("T16").Value = 0.03346139994672245

